$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("TC6_SearchResults_Typeahead")
$ws2 = $wb.Worksheets.Item("Testdata")

# Update the header text on sheet1
$ws1.Range("C8").Value = "ValidSearchHeader1"

# Set the (remembered) selection on sheet2 first, while it is still active,
# so selecting on it later doesn't re-activate it.
$ws2.Range("D9").Select()

# Make sheet1 the active sheet/tab, with C8 selected.
$ws1.Activate()
$ws1.Range("C8").Select()
